# Generate Report for Handback
# Update the generated timestamps that get refreshed each time the
# handback status report is produced.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the
# 83c2377b-...-a62b.md row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-06 22:56:46"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the same row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-06 22:56:41"
$wsZhCn.Range("K3").Value = "2016-09-06 22:57:23"

# de-de sheet: "Correspond Handoff Datetime" (shares text with the Overview
# sheet's value) / "Correspond Handback DateTime" for the same row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-06 22:56:46"
$wsDeDe.Range("K3").Value = "2016-09-06 22:57:31"
